$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3357.6936
$ws.Range("I76").Value = 3128.6726
$ws.Range("J76").Value = 5157.143
$ws.Range("K76").Value = 3128.6726
$ws.Range("L76").Value = 5157.143
$ws.Range("M76").Value = -2813.6726
$ws.Range("N76").Value = -5787.143
$ws.Range("H79").Value = 3357.6936
$ws.Range("I79").Value = 3128.6726
$ws.Range("J79").Value = 5157.143
$ws.Range("K79").Value = 3128.6726
$ws.Range("L79").Value = 5157.143
$ws.Range("M79").Value = -2036.6726
$ws.Range("N79").Value = -7341.143
$ws.Range("H112").Value = 1083.3
$ws.Range("I112").Value = 1125
$ws.Range("J112").Value = 1078.6666
$ws.Range("K112").Value = 3375
$ws.Range("L112").Value = 3235.9998
$ws.Range("M112").Value = -2267
$ws.Range("N112").Value = -5451.9998
$ws.Range("H121").Value = 1936.5
$ws.Range("I121").Value = 718
$ws.Range("J121").Value = 2490.3635
$ws.Range("K121").Value = 2154
$ws.Range("L121").Value = 7471.0905
$ws.Range("M121").Value = -407
$ws.Range("N121").Value = -10965.0905

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 16189.111
$ws.Range("I63").Value = 16189.111
$ws.Range("K63").Value = 16189.111
$ws.Range("M63").Value = -15503.111
$ws.Range("H66").Value = 16189.111
$ws.Range("I66").Value = 16189.111
$ws.Range("K66").Value = 80945.55500000001
$ws.Range("M66").Value = -77513.55500000001
$ws.Range("H132").Value = 2365.6875
$ws.Range("I132").Value = 1065.12
$ws.Range("J132").Value = 3779.348
$ws.Range("K132").Value = 3195.36
$ws.Range("L132").Value = 11338.044
$ws.Range("M132").Value = -665.3599999999997
$ws.Range("N132").Value = -16398.044

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1684.2858
$ws.Range("I105").Value = 1684.2858
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1684.2858
$ws.Range("L105").Value = 0
# row 105: M105 collapses to a single combined value and N105 is removed entirely
$ws.Range("M105").Value = 62.71419999999989
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 869.1177
$ws.Range("I107").Value = 899.36365
$ws.Range("J107").Value = 813.6667
$ws.Range("K107").Value = 899.36365
$ws.Range("L107").Value = 813.6667
$ws.Range("M107").Value = 1020.63635
$ws.Range("N107").Value = -4653.6667

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 38461620
$ws.Range("I12").Value = 166666670
$ws.Range("J12").Value = 107
$ws.Range("K12").Value = 500000010
$ws.Range("L12").Value = 321
$ws.Range("M12").Value = -499999837
$ws.Range("N12").Value = -667
$ws.Range("H64").Value = 989.5238000000001
$ws.Range("I64").Value = 500
$ws.Range("J64").Value = 1014
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 3042
$ws.Range("M64").Value = -1230
$ws.Range("N64").Value = -3582
$ws.Range("H67").Value = 989.5238000000001
$ws.Range("I67").Value = 500
$ws.Range("J67").Value = 1014
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 3042
$ws.Range("M67").Value = -564
$ws.Range("N67").Value = -4914
$ws.Range("H68").Value = 1699904.6
$ws.Range("J68").Value = 1868.6154
$ws.Range("L68").Value = 5605.8462
$ws.Range("N68").Value = -7227.8462
$ws.Range("H71").Value = 1699904.6
$ws.Range("J71").Value = 1868.6154
$ws.Range("L71").Value = 16817.5386
$ws.Range("N71").Value = -24929.5386
$ws.Range("H106").Value = 6864.5
$ws.Range("I106").Value = 4980
$ws.Range("J106").Value = 8749
$ws.Range("K106").Value = 14940
$ws.Range("L106").Value = 26247
$ws.Range("M106").Value = -13994
$ws.Range("N106").Value = -28139
$ws.Range("H109").Value = 2712.5715
$ws.Range("I109").Value = 763.3333
$ws.Range("J109").Value = 4174.5
$ws.Range("K109").Value = 2289.9999
$ws.Range("L109").Value = 12523.5
$ws.Range("M109").Value = -1249.9999
$ws.Range("N109").Value = -14603.5
$ws.Range("H123").Value = 1740.9231
$ws.Range("I123").Value = 839.8
$ws.Range("J123").Value = 2304.125
$ws.Range("K123").Value = 2519.4
$ws.Range("L123").Value = 6912.375
# row 123: M123 is a newly introduced cell between L123 and N123
$ws.Range("M123").Value = -69.39999999999964
$ws.Range("N123").Value = -11812.375
$ws.Range("H131").Value = 3208.638
$ws.Range("I131").Value = 531.5833
$ws.Range("J131").Value = 3907
$ws.Range("K131").Value = 1594.7499
$ws.Range("L131").Value = 11721
$ws.Range("M131").Value = 3445.2501
$ws.Range("N131").Value = -21801

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 55000
$ws.Range("J18").Value = 55000
$ws.Range("L18").Value = 55000
# row 18: N18 is a newly introduced cell (M18 remains absent)
$ws.Range("N18").Value = -55586
$ws.Range("H80").Value = 3737.1765
$ws.Range("I80").Value = 3804.5454
$ws.Range("J80").Value = 3613.6667
$ws.Range("K80").Value = 3804.5454
$ws.Range("L80").Value = 3613.6667
$ws.Range("M80").Value = -2806.5454
$ws.Range("N80").Value = -5609.6667
$ws.Range("H83").Value = 3737.1765
$ws.Range("I83").Value = 3804.5454
$ws.Range("J83").Value = 3613.6667
$ws.Range("K83").Value = 19022.727
$ws.Range("L83").Value = 18068.3335
$ws.Range("M83").Value = -14030.727
$ws.Range("N83").Value = -28052.3335
$ws.Range("H132").Value = 3317.7585
$ws.Range("I132").Value = 2134.647
$ws.Range("J132").Value = 4993.8335
$ws.Range("K132").Value = 6403.941
$ws.Range("L132").Value = 14981.5005
$ws.Range("M132").Value = -3873.941
$ws.Range("N132").Value = -20041.5005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45001.5
$ws.Range("I7").Value = 55938.26
$ws.Range("J7").Value = 3441.8
$ws.Range("K7").Value = 55938.26
$ws.Range("L7").Value = 3441.8
$ws.Range("M7").Value = -55826.26
$ws.Range("N7").Value = -3665.8
$ws.Range("H55").Value = 341.6111
$ws.Range("I55").Value = 167.81818
$ws.Range("J55").Value = 614.7143
$ws.Range("K55").Value = 167.81818
$ws.Range("L55").Value = 614.7143
$ws.Range("M55").Value = 5.181819999999988
$ws.Range("N55").Value = -960.7143
$ws.Range("H126").Value = 45001.5
$ws.Range("I126").Value = 55938.26
$ws.Range("J126").Value = 3441.8
$ws.Range("K126").Value = 167814.78
$ws.Range("L126").Value = 10325.4
$ws.Range("M126").Value = -165344.78
$ws.Range("N126").Value = -15265.4
$ws.Range("H136").Value = 5849729
$ws.Range("I136").Value = 1500.6342
$ws.Range("J136").Value = 20835814
$ws.Range("K136").Value = 4501.902599999999
$ws.Range("L136").Value = 62507442
$ws.Range("M136").Value = -1951.902599999999
$ws.Range("N136").Value = -62512542
